$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the print area from $A$1:$K$34 to $B$1:$K$34
$ws.PageSetup.PrintArea = '$B$1:$K$34'

# Update the selection from P28 to B1:K34
$ws.Range("B1:K34").Select()
